# Apply updated cryptos data (prices + 1h volume %) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) hold plain-text values (e.g. "65.540.64",
# "  -1.13%  ") rather than numbers. Force Text format first so Excel's COM
# layer does not auto-coerce numeric-looking strings ("600.02", "1.00", ...)
# into floating point numbers; restore the Normal style afterwards so no extra
# cell-level formatting lingers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "65.540.64"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "2.665.75"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "600.02"
$ws.Range("E5").Value = "  -1.94%  "
$ws.Range("D6").Value = "156.76"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.622"
$ws.Range("E8").Value = "  +5.26%  "
$ws.Range("D9").Value = "0.129"
$ws.Range("E9").Value = "  +2.44%  "
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("D11").Value = "5.85"
$ws.Range("E11").Value = "  -4.45%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "29.23"
$ws.Range("E13").Value = "  -4.09%  "
$ws.Range("E14").Value = "  -5.76%  "
$ws.Range("D15").Value = "3.145.53"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").Value = "65.393.15"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "2.671.35"
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "350.11"
$ws.Range("E21").Value = "  -3.61%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "69.62"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("D24").Value = "0.0000111"
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("E28").Value = "  -8.94%  "
$ws.Range("D29").Value = "8.09"
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "531.94"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").Value = "2.14"
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("E33").Value = "  -3.76%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "5.48"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").Value = "6.42"
$ws.Range("E35").Value = "  -4.94%  "
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "158.26"
$ws.Range("E39").Value = "  -3.30%  "
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("D42").Value = "42.83"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "164.28"
$ws.Range("E43").Value = "  -4.09%  "
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("D47").Value = "22.84"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("E49").Value = "  -3.70%  "
$ws.Range("D51").Value = "20.02"
$ws.Range("E51").Value = "  -5.28%  "

$ws.Range("D2:E51").Style = "Normal"
